$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5 (shifts current rows 5:25 down to 6:26)
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly price entry. The
# non-varying descriptive columns match the rest of the sheet (same
# market / product), and the week-specific values (date, volumes and
# prices) are the new observation.
$ws.Range("A5").Value2 = 1
$ws.Range("B5").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value2 = "Arica y Parinacota"
$ws.Range("D5").Value2 = 44882
$ws.Range("D5").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E5").Value2 = 15
$ws.Range("F5").Value2 = 100112003
$ws.Range("G5").Value2 = "Ajo"
$ws.Range("H5").Value2 = "Chino"
$ws.Range("I5").Value2 = "Primera"
$ws.Range("J5").Value2 = 400
$ws.Range("K5").Value2 = 15000
$ws.Range("L5").Value2 = 16000
$ws.Range("M5").Value2 = 15550
$ws.Range("N5").Value2 = "$/caja 10 kilos"
$ws.Range("O5").Value2 = "China"
$ws.Range("P5").Value2 = 1555
$ws.Range("Q5").Value2 = 10
$ws.Range("R5").Value2 = "Hortaliza"
